# Apply the "Added new signal read and baudrate to 9600" edit:
#  - rename pipeline_jmp -> pipeline_jmp_condl_rel_dests_cond_out (row 24),
#    and its bit-width label 1Bit -> 7Bit
#  - insert six new pipeline signal rows (26-31) describing the newly
#    exposed debug signals
#  - nudge the window/selection view state to match where the author was
#    working when they saved

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- row 24: pipeline_jmp became pipeline_jmp_condl_rel_dests_cond_out ---
$ws.Range("E24").Value = "pipeline_jmp_condl_rel_dests_cond_out"
$ws.Range("H24").Value = "7Bit"
$ws.Range("J24").Value = "pipeline_jmp_condl_rel_dests_cond_out"

# --- new rows 26-31: additional pipeline signals ---
$ws.Range("B26").Value = "0x18"
$ws.Range("E26").Value = "pipeline_immediate_out"
$ws.Range("G26").Value = "0x18"
$ws.Range("H26").Value = "16 Bit"
$ws.Range("J26").Value = "pipeline_immediate_out"

$ws.Range("B27").Value = "0x19"
$ws.Range("E27").Value = "pipeline_write_address_out"
$ws.Range("G27").Value = "0x19"
$ws.Range("H27").Value = "4 Bit"
$ws.Range("J27").Value = "pipeline_write_address_out"

$ws.Range("B28").Value = "0x1A"
$ws.Range("E28").Value = "pipeline_whb_wlb_out"
$ws.Range("G28").Value = "0x1A"
$ws.Range("H28").Value = "2 Bit"
$ws.Range("J28").Value = "pipeline_whb_wlb_out"

$ws.Range("B29").Value = "0x1B"
$ws.Range("E29").Value = "pipeline_write_data_sel_out"
$ws.Range("G29").Value = "0x1B"
$ws.Range("H29").Value = "2 Bit"
$ws.Range("J29").Value = "pipeline_write_data_sel_out"

$ws.Range("B30").Value = "0x1C"
$ws.Range("E30").Value = "pipeline_RAM_src_read_write_bankid_out"
$ws.Range("G30").Value = "0x1C"
$ws.Range("H30").Value = "7 Bit"
$ws.Range("J30").Value = "pipeline_RAM_src_read_write_bankid_out"

$ws.Range("B31").Value = "0x1D"
$ws.Range("E31").Value = "pipeline_is_alu_ram_gpu_op_out"
$ws.Range("G31").Value = "0x1D"
$ws.Range("H31").Value = "3 Bit"
$ws.Range("J31").Value = "pipeline_is_alu_ram_gpu_op_out"

# --- view/selection state as last saved by the author ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("D58").Select()
